$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add '.' to the degree/certification abbreviations in the Educational Background table
$ws.Range("A16").Value = "Ph.D. Boreal Ecology"
$ws.Range("A17").Value = "M.Sc. Biology"
$ws.Range("A18").Value = "B.Sc. Marine and Freshwater Biology"

# Update the last active selection to reflect where editing ended
$ws.Range("A18").Select()
